$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 175 ---------------------------------------------------------
# Column A (date): copy the date-formatted style from the last existing
# row (174) so the new cell reuses the same style index instead of
# Excel allocating a brand new cellXf.
$ws.Range("A174").Copy()
$ws.Range("A175").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A175").Value = 45485.2916666667

$ws.Range("B175").Value = 0
$ws.Range("C175").Value = 2.83999991416931
$ws.Range("D175").Value = 2.83999991416931
$ws.Range("E175").Value = 2.83999991416931
$ws.Range("F175").Value = 2.83999991416931

# Column G (adj_close, stored as text): the value is identical to an
# already-existing shared string (row 174's G cell), so copy that
# cell's value straight across - this reuses the existing shared
# string entry and keeps it typed as text without touching styles.xml.
$ws.Range("G174").Copy()
$ws.Range("G175").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("H175").Value = "XHS.MI"

# --- Row 176 ---------------------------------------------------------
$ws.Range("A174").Copy()
$ws.Range("A176").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A176").Value = 45488.3146527778

$ws.Range("B176").Value = 500
$ws.Range("C176").Value = 2.75999999046326
$ws.Range("D176").Value = 2.75999999046326
$ws.Range("E176").Value = 2.75999999046326
$ws.Range("F176").Value = 2.75999999046326

# Column G: this value is new, so force text formatting before
# assigning the numeric-looking string (otherwise Excel would parse it
# back into a number), then drop back to the Normal style so no
# lingering number format is left attached to the cell itself.
$ws.Range("G176").NumberFormat = "@"
$ws.Range("G176").Value = "2.75999999046326"
$ws.Range("G176").Style = "Normal"

$ws.Range("H176").Value = "XHS.MI"

$excel.CutCopyMode = 0
